$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 464.5
$ws.Range("I4").Value = 286
$ws.Range("K4").Value = 286
$ws.Range("M4").Value = -172
$ws.Range("H21").Value = 23846.154
$ws.Range("J21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("N21").Value = -25936
$ws.Range("H23").Value = 23846.154
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25468
$ws.Range("H29").Value = 400
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 10352.632
$ws.Range("I62").Value = 10394.235
$ws.Range("K62").Value = 10394.235
$ws.Range("M62").Value = -9770.235000000001
$ws.Range("H65").Value = 10352.632
$ws.Range("I65").Value = 10394.235
$ws.Range("K65").Value = 51971.175
$ws.Range("M65").Value = -48851.175
$ws.Range("H86").Value = 1993
$ws.Range("I86").Value = 1991
$ws.Range("J86").Value = 1995
$ws.Range("K86").Value = 1991
$ws.Range("L86").Value = 1995
$ws.Range("M86").Value = -868
$ws.Range("N86").Value = -4241
$ws.Range("H89").Value = 1993
$ws.Range("I89").Value = 1991
$ws.Range("J89").Value = 1995
$ws.Range("K89").Value = 9955
$ws.Range("L89").Value = 9975
$ws.Range("M89").Value = -4339
$ws.Range("N89").Value = -21207
$ws.Range("H112").Value = 3437
$ws.Range("J112").Value = 3437
$ws.Range("L112").Value = 10311
$ws.Range("N112").Value = -12527
$ws.Range("H138").Value = 2498.2222
$ws.Range("I138").Value = 2247
$ws.Range("J138").Value = 2760.8635
$ws.Range("K138").Value = 6741
$ws.Range("L138").Value = 8282.5905
$ws.Range("M138").Value = -1601
$ws.Range("N138").Value = -18562.5905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2639.75
$ws.Range("I94").Value = 2025.579
$ws.Range("J94").Value = 3936.3333
$ws.Range("K94").Value = 2025.579
$ws.Range("L94").Value = 3936.3333
$ws.Range("M94").Value = -1574.579
$ws.Range("N94").Value = -4838.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2848.8
$ws.Range("I31").Value = 1873.5
$ws.Range("J31").Value = 6750
$ws.Range("K31").Value = 1873.5
$ws.Range("L31").Value = 6750
$ws.Range("M31").Value = -1578.5
$ws.Range("N31").Value = -7340
$ws.Range("H34").Value = 2848.8
$ws.Range("I34").Value = 1873.5
$ws.Range("J34").Value = 6750
$ws.Range("K34").Value = 1873.5
$ws.Range("L34").Value = 6750
$ws.Range("M34").Value = -1671.5
$ws.Range("N34").Value = -7154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 596.125
$ws.Range("I5").Value = 610
$ws.Range("K5").Value = 1830
$ws.Range("M5").Value = -1718
$ws.Range("H50").Value = 125206.875
$ws.Range("I50").Value = 185
$ws.Range("J50").Value = 166880.83
$ws.Range("K50").Value = 555
$ws.Range("L50").Value = 500642.49
$ws.Range("M50").Value = -74
$ws.Range("N50").Value = -501604.49
$ws.Range("H53").Value = 125206.875
$ws.Range("I53").Value = 185
$ws.Range("J53").Value = 166880.83
$ws.Range("K53").Value = 555
$ws.Range("L53").Value = 500642.49
$ws.Range("M53").Value = -74
$ws.Range("N53").Value = -501604.49
$ws.Range("H69").Value = 949.5
$ws.Range("I69").Value = 949.5
$ws.Range("K69").Value = 2848.5
$ws.Range("M69").Value = -2037.5
$ws.Range("H72").Value = 949.5
$ws.Range("I72").Value = 949.5
$ws.Range("K72").Value = 8545.5
$ws.Range("M72").Value = -4489.5
$ws.Range("H135").Value = 596.125
$ws.Range("I135").Value = 610
$ws.Range("K135").Value = 5490
$ws.Range("M135").Value = -2955

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7469.8
$ws.Range("I70").Value = 6849.857
$ws.Range("K70").Value = 6849.857
$ws.Range("M70").Value = -6579.857
$ws.Range("H73").Value = 7469.8
$ws.Range("I73").Value = 6849.857
$ws.Range("K73").Value = 6849.857
$ws.Range("M73").Value = -5913.857
$ws.Range("H132").Value = 33176.426
$ws.Range("I132").Value = 58495.277
$ws.Range("J132").Value = 2793.8
$ws.Range("K132").Value = 175485.831
$ws.Range("L132").Value = 8381.400000000001
$ws.Range("M132").Value = -172955.831
$ws.Range("N132").Value = -13441.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2500.3333
$ws.Range("I7").Value = 2364
$ws.Range("K7").Value = 2364
$ws.Range("M7").Value = -2252
$ws.Range("H40").Value = 7251.4116
$ws.Range("I40").Value = 6981.4614
$ws.Range("J40").Value = 8128.75
$ws.Range("K40").Value = 6981.4614
$ws.Range("L40").Value = 8128.75
$ws.Range("M40").Value = -6845.4614
$ws.Range("N40").Value = -8400.75
$ws.Range("H44").Value = 29999.5
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10912
$ws.Range("H68").Value = 5874.25
$ws.Range("I68").Value = 6949
$ws.Range("J68").Value = 4799.5
$ws.Range("K68").Value = 6949
$ws.Range("L68").Value = 4799.5
$ws.Range("M68").Value = -6200
$ws.Range("N68").Value = -6297.5
$ws.Range("H71").Value = 5874.25
$ws.Range("I71").Value = 6949
$ws.Range("J71").Value = 4799.5
$ws.Range("K71").Value = 34745
$ws.Range("L71").Value = 23997.5
$ws.Range("M71").Value = -31001
$ws.Range("N71").Value = -31485.5
$ws.Range("H104").Value = 16061
$ws.Range("J104").Value = 16061
$ws.Range("L104").Value = 16061
$ws.Range("N104").Value = -23049
$ws.Range("H126").Value = 2500.3333
$ws.Range("I126").Value = 2364
$ws.Range("K126").Value = 7092
$ws.Range("M126").Value = -4622
$ws.Range("H127").Value = 87875
$ws.Range("J127").Value = 87875
$ws.Range("L127").Value = 87875
$ws.Range("N127").Value = -97795

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2818.8333
$ws.Range("I122").Value = 2230.7778
$ws.Range("K122").Value = 6692.3334
$ws.Range("M122").Value = -4242.3334
$ws.Range("H136").Value = 3371.6
$ws.Range("I136").Value = 3403.913
$ws.Range("K136").Value = 10211.739
$ws.Range("M136").Value = -7661.739
